# PM04 Tidsregistrering for Bille.xlsx - add new time-registration entries
# for rows 24-29 (previously blank template rows) and move the active
# selection, matching the author's newest timesheet entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 24: UI (FX) for Vareforbrug / User-Interface Designer / 05-03-2020
$ws.Range("A23:E23").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)
$ws.Range("A24").Value = "UI (FX) for Vareforbrug"
$ws.Range("B24").Value = "User-Interface Designer"
$ws.Range("C24").Value = 43895
$ws.Range("D24").Value = 0.375
$ws.Range("E24").Value = 0.45833333333333331

# Row 25: UC0804 / Implenter / 05-03-2020
$ws.Range("A23:E23").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A25").Value = "UC0804"
$ws.Range("B25").Value = "Implenter"
$ws.Range("C25").Value = 43895
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = 0.61458333333333337

# Row 26: sammenlæg dd og dd02 / (no role) / 06-03-2020
$ws.Range("A23").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("C23:E23").Copy()
$ws.Range("C26:E26").PasteSpecial(-4122)
$ws.Range("A26").Value = "sammenlæg dd og dd02"
$ws.Range("C26").Value = 43896
$ws.Range("D26").Value = 0.39583333333333331
$ws.Range("E26").Value = 0.41666666666666669

# Row 27: STD0104 / Test Analyst / 06-03-2020
$ws.Range("A23:E23").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)
$ws.Range("A27").Value = "STD0104"
$ws.Range("B27").Value = "Test Analyst"
$ws.Range("C27").Value = 43896
$ws.Range("D27").Value = 0.46875
$ws.Range("E27").Value = 0.53125

# Row 28: UnitTest- STD0104 / Test Desinger / 06-03-2020
$ws.Range("A23:E23").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)
$ws.Range("A28").Value = "UnitTest- STD0104"
$ws.Range("B28").Value = "Test Desinger"
$ws.Range("C28").Value = 43896
$ws.Range("D28").Value = 0.53125
$ws.Range("E28").Value = 0.60416666666666663

# Row 29: UnitTest- STD0104 / Test Desinger / 09-03-2020
$ws.Range("A23:E23").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Range("A29").Value = "UnitTest- STD0104"
$ws.Range("B29").Value = "Test Desinger"
$ws.Range("C29").Value = 43899
$ws.Range("D29").Value = 0.375
$ws.Range("E29").Value = 0.54166666666666663

# Move the viewport / active selection to where the user last was.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D30").Select() | Out-Null
